$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated cryptocurrency price/volume snapshot values (Price column D, Volume(1h) column E).
# Target cells are stored as text (e.g. "328.63", "-0.29%"), so force text format
# before assignment to avoid Excel re-interpreting them as numbers/percentages.
$updates = @{
    "D2" = "328.63"
    "E2" = "-0.29%"
    "D3" = "44.36"
    "E3" = "0.14%"
    "D4" = "5.029"
    "E4" = "-8.75%"
    "D5" = "0.08389"
    "E5" = "3.95%"
    "D6" = "1.951"
    "E6" = "-6.09%"
    "D7" = "0.9775"
    "E7" = "1.69%"
    "E8" = "-4.47%"
    "D9" = "0.1167"
    "E9" = "2.94%"
    "D10" = "0.1897"
    "E10" = "0.82%"
    "D11" = "0.09693"
    "E11" = "-2.44%"
    "D12" = "0.04629"
    "E12" = "-1.14%"
    "D13" = "0.1060"
    "E13" = "0.43%"
    "D14" = "0.001292"
    "E14" = "1.89%"
    "D15" = "0.005919"
    "E15" = "-2.88%"
    "D16" = "3.404"
    "E16" = "1.94%"
    "D17" = "4.439"
    "E17" = "0.04%"
    "D18" = "0.3321"
    "E18" = "0.29%"
    "D19" = "8.915"
    "E19" = "-12.34%"
    "D20" = "0.1352"
    "E20" = "-2.41%"
    "D21" = "0.2551"
    "E21" = "-1.07%"
    "D22" = "0.04165"
    "E22" = "1.52%"
    "D23" = "0.001299"
    "E23" = "-0.92%"
    "D24" = "0.004545"
    "E24" = "3.88%"
    "D25" = "0.0001302"
    "E25" = "1.59%"
    "D26" = "0.0002982"
    "E26" = "-20.42%"
    "D38" = "0.02739"
    "E38" = "2.85%"
    "D39" = "0.05651"
    "E39" = "0.18%"
    "D40" = "0.007861"
    "E40" = "2.82%"
    "D41" = "0.1417"
    "E41" = "0.29%"
    "D42" = "0.007361"
    "E42" = "-0.33%"
    "D43" = "0.002070"
    "E43" = "4.09%"
    "D44" = "0.007902"
    "E44" = "-9.18%"
    "D45" = "0.3509"
    "D46" = "0.00006911"
    "E46" = "-2.95%"
    "D47" = "0.00000000751"
    "E47" = "-0.01%"
    "D48" = "0.003512"
    "E48" = "1.89%"
    "D49" = "0.003534"
    "E49" = "39.94%"
    "D50" = "0.00002104"
    "E50" = "-0.01%"
    "D51" = "0.0002003"
    "E51" = "-0.01%"
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$addr]
}
